# Update the dSF ("F") column values for a handful of rows on the active sheet.
# These correspond to a "repull data / push all data / mean calculation" update
# where the final dS (dSF) values were recalculated for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -1
$ws.Range("F5").Value  = -2
$ws.Range("F10").Value = -12
$ws.Range("F24").Value = -1
$ws.Range("F26").Value = 1
$ws.Range("F42").Value = -1
$ws.Range("F50").Value = 1
$ws.Range("F51").Value = -3
